# #12 Issues synthesis added on slides
#
# Insert a new "Title and Content" slide at show-position 11 (right after the
# "SAI" slide / before the "Solutions" slide) containing an "Issues" title
# and a two-line "Issues synthesis" / "Use a complex shema" body.

$p = $ppt.ActivePresentation

# ppLayoutText (2) == the "Title and Content" autolayout used by every other
# content slide in this deck (slideLayout2.xml, "Titre et contenu").
$s = $p.Slides.Add(11, 2)

# Title placeholder
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Issues"

# Body / content placeholder
$body = $s.Shapes.Item(2)
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "Issues synthesis" + [char]13 + "Use a complex shema"
